# Apply updated crypto price/volume figures (data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.038.77'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '1.822.27'
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = '''312.04'
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").Value = '''0.4560'
$ws.Range("E7").Value = '  +6.83%  '
$ws.Range("D8").Value = '''0.3712'
$ws.Range("E8").Value = '  +1.15%  '
$ws.Range("D9").Value = '''0.07291'
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("D10").Value = '''0.8581'
$ws.Range("E10").Value = '  -0.29%  '
$ws.Range("D11").Value = '''20.82'
$ws.Range("E11").Value = '  -0.79%  '
$ws.Range("D12").Value = '1.824.88'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = '''6.656'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").Value = '''92.86'
$ws.Range("E14").Value = '  +5.27%  '
$ws.Range("D15").Value = '''5.332'
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").Value = '''0.07093'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("D18").Value = '''0.000008822'
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").Value = '27.027.04'
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("D22").Value = '''5.175'
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("D23").Value = '''10.94'
$ws.Range("E23").Value = '  +0.83%  '
$ws.Range("D24").Value = '''1.992'
$ws.Range("E24").Value = '  -0.54%  '
$ws.Range("D25").Value = '''151.62'
$ws.Range("E25").Value = '  -1.25%  '
$ws.Range("D26").Value = '''2.220'
$ws.Range("E26").Value = '  +4.91%  '
$ws.Range("D27").Value = '''18.43'
$ws.Range("E27").Value = '  +0.72%  '
$ws.Range("D28").Value = '''5.255'
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("D29").Value = '''116.78'
$ws.Range("E29").Value = '  +0.57%  '
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("D32").Value = '''0.7533'
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("D33").Value = '''2.939'
$ws.Range("E33").Value = '  +3.89%  '
$ws.Range("D34").Value = '''4.457'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").Value = '''1.001'
$ws.Range("E35").Value = '  -0.43%  '
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("D37").Value = '''0.01966'
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("D38").Value = '''0.05250'
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("D39").Value = '''0.5326'
$ws.Range("E39").Value = '  +5.62%  '
$ws.Range("D40").Value = '''7.194'
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("D41").Value = '''2.877'
$ws.Range("E41").Value = '  -0.79%  '
$ws.Range("D42").Value = '''0.1716'
$ws.Range("E42").Value = '  +1.86%  '
$ws.Range("D43").Value = '''0.5208'
$ws.Range("E43").Value = '  +9.72%  '
$ws.Range("D44").Value = '''8.543'
$ws.Range("E44").Value = '  -0.66%  '
$ws.Range("D45").Value = '''10.65'
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("D46").Value = '''1.967'
$ws.Range("E46").Value = '  +9.14%  '
$ws.Range("D47").Value = '''105.55'
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("D48").Value = '''1.672'
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("D49").Value = '''1.000'
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("D50").Value = '''0.06385'
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").Value = '''63.31'
$ws.Range("E51").Value = '  +0.31%  '
